$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31
$ws.Range("A31").Value = "Philippines"
$ws.Range("B31").Value = "philippines"
$ws.Range("C31").Value = "Bicol River Basin"
$ws.Range("D31").Value = "bicol"
$ws.Range("E31").Value = "Nabua"
$ws.Range("F31").Value = "G4611"
$ws.Range("G31").Value = "primary"
$ws.Range("H31").Value = "'2025-10-30"
$ws.Range("I31").Value = 3
$ws.Range("J31").Value = 13.37499999999993
$ws.Range("K31").Value = 123.3249999999996
$ws.Range("L31").Value = 5
$ws.Range("M31").Value = 483.4323679605675
$ws.Range("N31").Value = "LOW"
$ws.Range("O31").Value = 348.6773053168241
$ws.Range("P31").Value = 483.4323679605675
$ws.Range("Q31").Value = 50
$ws.Range("R31").Value = 0
$ws.Range("S31").Value = 0
$ws.Range("T31").Value = 66.09375
$ws.Range("U31").Value = 67.83062744140625
$ws.Range("V31").Value = 56.0078125
$ws.Range("W31").Value = 97.09375
$ws.Range("X31").Value = 60.806640625
$ws.Range("Y31").Value = 72.181640625
$ws.Range("Z31").Value = $false
$ws.Range("AA31").Value = -86.32823236912611

# Row 32
$ws.Range("A32").Value = "Philippines"
$ws.Range("B32").Value = "philippines"
$ws.Range("C32").Value = "Bicol River Basin"
$ws.Range("D32").Value = "bicol"
$ws.Range("E32").Value = "Nabua"
$ws.Range("F32").Value = "G4611"
$ws.Range("G32").Value = "primary"
$ws.Range("H32").Value = "'2025-10-31"
$ws.Range("I32").Value = 3
$ws.Range("J32").Value = 13.37499999999993
$ws.Range("K32").Value = 123.3249999999996
$ws.Range("L32").Value = 5
$ws.Range("M32").Value = 483.4323679605675
$ws.Range("N32").Value = "LOW"
$ws.Range("O32").Value = 348.6773053168241
$ws.Range("P32").Value = 483.4323679605675
$ws.Range("Q32").Value = 50
$ws.Range("R32").Value = 0
$ws.Range("S32").Value = 0
$ws.Range("T32").Value = 40.0625
$ws.Range("U32").Value = 41.31937408447266
$ws.Range("V32").Value = 33.265625
$ws.Range("W32").Value = 66.4609375
$ws.Range("X32").Value = 37.783203125
$ws.Range("Y32").Value = 43.517578125
$ws.Range("Z32").Value = $false
$ws.Range("AA32").Value = -91.7129049159431

# Strip the quote-prefix style picked up when forcing date-like strings to text,
# so the new cells keep the workbook default style (matches the source rows).
$ws.Range("H31").ClearFormats()
$ws.Range("H32").ClearFormats()
